# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E) and Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

# Per-row updates to columns C (index), E (cumulative value) and G (pulse width).
# D (Point_Exceeds_Index) and F (Point_Exceeds_Cumulative_Value) are unchanged.
$rowUpdates = @{
    2 = @{ C = 47; E = 0.001086886039829613 }
    3 = @{ C = 50; E = 0.03227880228351791 }
    5 = @{ C = 47; E = 0.0006534875571948642 }
    6 = @{ C = 50; E = 0.05006205891674615 }
}

# Pulse_Width (column G) values differ per sheet.
$gUpdates = @{
    "Step3_DataPts_0.5" = @{ 2 = 51; 3 = 48; 5 = 49; 6 = 47 }
    "Step3_DataPts_0.7" = @{ 2 = 72; 3 = 70; 5 = 62; 6 = 63 }
    "Step3_DataPts_0.8" = @{ 2 = 86; 3 = 85; 5 = 79; 6 = 76 }
    "Step3_DataPts_0.9" = @{ 2 = 120; 3 = 118; 5 = 114; 6 = 117 }
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    foreach ($row in $rowUpdates.Keys) {
        $vals = $rowUpdates[$row]
        $ws.Cells.Item($row, 3).Value = $vals.C   # Column C
        $ws.Cells.Item($row, 5).Value = $vals.E   # Column E
    }

    $gForSheet = $gUpdates[$name]
    foreach ($row in $gForSheet.Keys) {
        $ws.Cells.Item($row, 7).Value = $gForSheet[$row]  # Column G
    }
}
